# Auto-generated script applying scheduled market-price refresh values
# to the profit-calculation columns (H-N) across all class sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2500889.2
$ws.Range("J17").Value = 2500889.2
$ws.Range("L17").Value = 7502667.600000001
$ws.Range("N17").Value = -7503003.600000001

$ws.Range("H74").Value = 4625.75
$ws.Range("I74").Value = 3503
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 3503
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -2567
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 4625.75
$ws.Range("I77").Value = 3503
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 17515
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -12835
$ws.Range("N77").Value = -34360

$ws.Range("H98").Value = 991.6286
$ws.Range("I98").Value = 995.1
$ws.Range("J98").Value = 987
$ws.Range("K98").Value = 995.1
$ws.Range("L98").Value = 987
$ws.Range("M98").Value = 502.9
$ws.Range("N98").Value = -3983

$ws.Range("H112").Value = 1214.7069
$ws.Range("I112").Value = 736
$ws.Range("J112").Value = 1259.8679
$ws.Range("K112").Value = 2208
$ws.Range("L112").Value = 3779.6037
$ws.Range("M112").Value = -1100
$ws.Range("N112").Value = -5995.6037

$ws.Range("H122").Value = 991.6286
$ws.Range("I122").Value = 995.1
$ws.Range("J122").Value = 987
$ws.Range("K122").Value = 2985.3
$ws.Range("L122").Value = 2961
$ws.Range("M122").Value = -535.3000000000002
$ws.Range("N122").Value = -7861

$ws.Range("H132").Value = 1282.9841
$ws.Range("I132").Value = 576.62964
$ws.Range("J132").Value = 5521.1113
$ws.Range("K132").Value = 1729.88892
$ws.Range("L132").Value = 16563.3339
$ws.Range("M132").Value = 800.1110800000001
$ws.Range("N132").Value = -21623.3339

$ws.Range("H137").Value = 1891.0541
$ws.Range("I137").Value = 1529.579
$ws.Range("J137").Value = 2272.611
$ws.Range("K137").Value = 4588.737
$ws.Range("L137").Value = 6817.833
$ws.Range("M137").Value = -2038.737
$ws.Range("N137").Value = -11917.833

$ws.Range("H138").Value = 2251.6956
$ws.Range("I138").Value = 1624.7368
$ws.Range("J138").Value = 3020.2258
$ws.Range("K138").Value = 4874.2104
$ws.Range("L138").Value = 9060.6774
$ws.Range("M138").Value = 265.7896000000001
$ws.Range("N138").Value = -19340.6774

$ws.Range("H141").Value = 6056.0527
$ws.Range("I141").Value = 2275.3572
$ws.Range("J141").Value = 16642
$ws.Range("K141").Value = 6826.071599999999
$ws.Range("L141").Value = 49926
$ws.Range("M141").Value = -1646.071599999999
$ws.Range("N141").Value = -60286

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7569.949
$ws.Range("I32").Value = 5959.9185
$ws.Range("K32").Value = 5959.9185
$ws.Range("M32").Value = -5672.9185

$ws.Range("H45").Value = 1307.303
$ws.Range("I45").Value = 1263.963
$ws.Range("J45").Value = 1502.3334
$ws.Range("K45").Value = 1263.963
$ws.Range("L45").Value = 1502.3334
$ws.Range("M45").Value = -886.963
$ws.Range("N45").Value = -2256.3334

$ws.Range("H61").Value = 3756.7778
$ws.Range("I61").Value = 7706
$ws.Range("J61").Value = 2628.4285
$ws.Range("K61").Value = 7706
$ws.Range("L61").Value = 2628.4285
$ws.Range("M61").Value = -7494
$ws.Range("N61").Value = -3052.4285

$ws.Range("H132").Value = 1551.7906
$ws.Range("I132").Value = 938.2646999999999
$ws.Range("K132").Value = 2814.7941
$ws.Range("M132").Value = -284.7941000000001

$ws.Range("H136").Value = 3756.7778
$ws.Range("I136").Value = 7706
$ws.Range("J136").Value = 2628.4285
$ws.Range("K136").Value = 23118
$ws.Range("L136").Value = 7885.2855
$ws.Range("M136").Value = -20568
$ws.Range("N136").Value = -12985.2855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1617.6923
$ws.Range("I99").Value = 995
$ws.Range("J99").Value = 1730.909
$ws.Range("K99").Value = 995
$ws.Range("L99").Value = 1730.909
$ws.Range("M99").Value = 503
$ws.Range("N99").Value = -4726.909

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49974.25
$ws.Range("J20").Value = 49974.25
$ws.Range("L20").Value = 49974.25
$ws.Range("N20").Value = -50446.25

$ws.Range("H30").Value = 49974.25
$ws.Range("J30").Value = 49974.25
$ws.Range("L30").Value = 49974.25
$ws.Range("N30").Value = -50156.25

$ws.Range("H128").Value = 49974.25
$ws.Range("J128").Value = 49974.25
$ws.Range("L128").Value = 49974.25
$ws.Range("N128").Value = -59934.25

$ws.Range("H132").Value = 2229.76
$ws.Range("I132").Value = 1692.1177
$ws.Range("J132").Value = 3372.25
$ws.Range("K132").Value = 5076.3531
$ws.Range("L132").Value = 10116.75
$ws.Range("M132").Value = -2546.3531
$ws.Range("N132").Value = -15176.75

$ws.Range("H134").Value = 2483.818
$ws.Range("I134").Value = 2569.111
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 7707.333
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -5172.333
$ws.Range("N134").Value = -11370

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3737.375
$ws.Range("I3").Value = 2316.5
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 6949.5
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = -6837.5
$ws.Range("N3").Value = -24224

$ws.Range("H5").Value = 1254277.1
$ws.Range("I5").Value = 426.875
$ws.Range("J5").Value = 2926077.2
$ws.Range("K5").Value = 1280.625
$ws.Range("L5").Value = 8778231.600000001
$ws.Range("M5").Value = -1168.625
$ws.Range("N5").Value = -8778455.600000001

$ws.Range("H131").Value = 858.9299999999999
$ws.Range("I131").Value = 416
$ws.Range("J131").Value = 882.2421000000001
$ws.Range("K131").Value = 1248
$ws.Range("L131").Value = 2646.7263
$ws.Range("M131").Value = 3792
$ws.Range("N131").Value = -12726.7263

$ws.Range("H133").Value = 1898
$ws.Range("I133").Value = 1796
$ws.Range("K133").Value = 5388
$ws.Range("M133").Value = -328

$ws.Range("H135").Value = 1254277.1
$ws.Range("I135").Value = 426.875
$ws.Range("J135").Value = 2926077.2
$ws.Range("K135").Value = 3841.875
$ws.Range("L135").Value = 26334694.8
$ws.Range("M135").Value = -1306.875
$ws.Range("N135").Value = -26339764.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1717.8334
$ws.Range("I2").Value = 2582.5
$ws.Range("J2").Value = 1285.5
$ws.Range("K2").Value = 2582.5
$ws.Range("L2").Value = 1285.5
$ws.Range("M2").Value = -2469.5
$ws.Range("N2").Value = -1511.5

$ws.Range("H97").Value = 1037.4445
$ws.Range("I97").Value = 1094.2667
$ws.Range("J97").Value = 753.3333
$ws.Range("K97").Value = 1094.2667
$ws.Range("L97").Value = 753.3333
$ws.Range("M97").Value = -598.2666999999999
$ws.Range("N97").Value = -1745.3333

$ws.Range("H126").Value = 2200.6875
$ws.Range("I126").Value = 1907.2222
$ws.Range("J126").Value = 2578
$ws.Range("K126").Value = 5721.6666
$ws.Range("L126").Value = 7734
$ws.Range("M126").Value = -3251.6666
$ws.Range("N126").Value = -12674

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12381.552
$ws.Range("I132").Value = 15062.667
$ws.Range("J132").Value = 7994.273
$ws.Range("K132").Value = 45188.001
$ws.Range("L132").Value = 23982.819
$ws.Range("M132").Value = -42658.001
$ws.Range("N132").Value = -29042.819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1128
$ws.Range("I113").Value = 1176
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 3528
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -1358
$ws.Range("N113").Value = -6140

$ws.Range("H128").Value = 49424.375
$ws.Range("J128").Value = 49424.375
$ws.Range("L128").Value = 49424.375
$ws.Range("N128").Value = -59384.375

$ws.Range("H138").Value = 48411.4
$ws.Range("J138").Value = 48411.4
$ws.Range("L138").Value = 48411.4
$ws.Range("N138").Value = -58691.4
